# Update cryptocurrency price (column D) and volume change (column E) values
# to match the latest scrape, as produced by the scheduled GitHub Actions run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    # Prefix with an apostrophe so Excel always stores the value as text,
    # even when it looks numeric (e.g. '0.9989' or '26.389.54'), then reset
    # the cell style back to Normal so no quote-prefix / text format sticks.
    $cell = $ws.Range($cellRef)
    $cell.Value2 = "'" + $text
    $cell.Style = "Normal"
}

Set-TextValue "D2" '26.389.54'
Set-TextValue "E2" '  -0.37%  '
Set-TextValue "D3" '1.724.52'
Set-TextValue "E3" '  -0.43%  '
Set-TextValue "D4" '0.9989'
Set-TextValue "D5" '242.62'
Set-TextValue "E5" '  -1.94%  '
Set-TextValue "D7" '0.4880'
Set-TextValue "E7" '  +0.07%  '
Set-TextValue "D8" '0.2591'
Set-TextValue "E8" '  -2.99%  '
Set-TextValue "D9" '0.06193'
Set-TextValue "E9" '  -0.55%  '
Set-TextValue "D10" '1.727.61'
Set-TextValue "E10" '  -0.26%  '
Set-TextValue "D11" '0.06981'
Set-TextValue "E11" '  -1.20%  '
Set-TextValue "D12" '15.53'
Set-TextValue "E12" '  -0.82%  '
Set-TextValue "D13" '4.525'
Set-TextValue "E13" '  -2.90%  '
Set-TextValue "D14" '0.5979'
Set-TextValue "E14" '  -1.95%  '
Set-TextValue "D15" '77.19'
Set-TextValue "E15" '  -0.36%  '
Set-TextValue "D16" '0.9992'
Set-TextValue "E16" '  -0.09%  '
Set-TextValue "D17" '26.393.38'
Set-TextValue "E17" '  -0.38%  '
Set-TextValue "D18" '0.9992'
Set-TextValue "E18" '  -0.08%  '
Set-TextValue "D19" '0.000007200'
Set-TextValue "E19" '  +0.43%  '
Set-TextValue "E20" '  -1.68%  '
Set-TextValue "D21" '1.940.60'
Set-TextValue "E21" '  -1.02%  '
Set-TextValue "D22" '4.448'
Set-TextValue "E22" '  -1.76%  '
Set-TextValue "D23" '8.496'
Set-TextValue "E23" '  -3.25%  '
Set-TextValue "D24" '5.104'
Set-TextValue "E24" '  -2.98%  '
Set-TextValue "D25" '138.16'
Set-TextValue "E25" '  -0.82%  '
Set-TextValue "E26" '  -1.23%  '
Set-TextValue "E27" '  -0.24%  '
Set-TextValue "D28" '106.30'
Set-TextValue "E28" '  -1.67%  '
Set-TextValue "D29" '1.727'
Set-TextValue "E29" '  -2.79%  '
Set-TextValue "E30" '  -1.70%  '
Set-TextValue "D31" '0.08012'
Set-TextValue "E31" '  -0.29%  '
Set-TextValue "D32" '3.656'
Set-TextValue "E32" '  -1.19%  '
Set-TextValue "D33" '0.04500'
Set-TextValue "E33" '  -1.80%  '
Set-TextValue "D34" '2.603'
Set-TextValue "E34" '  -0.43%  '
Set-TextValue "D35" '0.9988'
Set-TextValue "E35" '  -0.53%  '
Set-TextValue "D36" '0.6235'
Set-TextValue "E36" '  -2.11%  '
Set-TextValue "D37" '0.9292'
Set-TextValue "E37" '  +3.63%  '
Set-TextValue "D38" '1.961'
Set-TextValue "E38" '  -2.88%  '
Set-TextValue "D39" '2.379'
Set-TextValue "E39" '  -0.80%  '
Set-TextValue "E40" '  -0.32%  '
Set-TextValue "D41" '0.01474'
Set-TextValue "E41" '  -2.31%  '
Set-TextValue "D42" '100.52'
Set-TextValue "E42" '  -0.95%  '
Set-TextValue "D43" '5.468'
Set-TextValue "E43" '  +0.39%  '
Set-TextValue "D44" '0.3843'
Set-TextValue "E44" '  -1.23%  '
Set-TextValue "D45" '6.922'
Set-TextValue "E45" '  -0.63%  '
Set-TextValue "D46" '0.1165'
Set-TextValue "E46" '  -1.60%  '
Set-TextValue "D47" '0.05362'
Set-TextValue "E47" '  -0.43%  '
Set-TextValue "D48" '30.17'
Set-TextValue "E48" '  -1.49%  '
Set-TextValue "D49" '7.713'
Set-TextValue "E49" '  -1.43%  '
Set-TextValue "D50" '1.227'
Set-TextValue "E50" '  -1.79%  '
Set-TextValue "D51" '50.98'
Set-TextValue "E51" '  -0.77%  '
